$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values in column C (Max column)
$ws.Range("C2").Value = 11.7
$ws.Range("C3").Value = 10.4

# Add new cells in columns E and F for rows 2 and 3, carrying over the
# same cell style used by the other data cells (copy format from B2,
# which already uses that style), then clear their value so they stay
# visually empty placeholders.
$ws.Range("B2").Copy()
$ws.Range("E2:F3").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("E2:F3").Value = ""

# Reflect the new selection state (active cell E2, selection E2:F3)
$ws.Range("E2:F3").Select()
